$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Copy existing row 7 (Miranda, older submission) down into new row 8 ---
# Row 8 keeps the original text-formatted SmartScore values (with trailing zeros)
# and gets its own timestamp for the new submission.

$ws.Cells.Item(8,1).Value = "Miranda"
$ws.Cells.Item(8,2).Value = 25
$ws.Cells.Item(8,3).Value = "Femenino"
$ws.Cells.Item(8,4).Value = "2025-10-28 06:02:28"
$ws.Cells.Item(8,5).Value = "{
  `"portion`": 0.8,
  `"diet`": 0.5714285714285714,
  `"salt`": 0.6,
  `"fat`": 0.8,
  `"natural`": 0.6,
  `"convenience`": 0.4,
  `"price`": 0.8
}"
$ws.Cells.Item(8,6).Value = "Nongshim Neoguri Spicy Seafood"
$ws.Cells.Item(8,8).Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Cells.Item(8,9).Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Cells.Item(8,11).Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Cells.Item(8,12).Value = "Maruchan Ramen Sabor Pollo"
$ws.Cells.Item(8,14).Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Cells.Item(8,15).Value = "Kraft Macaroni & Cheese Dinner"
$ws.Cells.Item(8,17).Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Cells.Item(8,18).Value = "Annie’s Shells & White Cheddar"
$ws.Cells.Item(8,20).Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Cells.Item(8,21).Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Cells.Item(8,23).Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Cells.Item(8,24).Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Cells.Item(8,26).Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Cells.Item(8,27).Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Cells.Item(8,29).Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Cells.Item(8,30).Value = "Jack Link’s Beef Jerky Original"
$ws.Cells.Item(8,32).Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"

# The multi-line JSON cell (E8) triggers an implicit autofit taller than default;
# settle the row height back to standard BEFORE touching per-cell number formats
# below (AutoFit after a NumberFormat/ClearFormats dance can resurrect the style).
$ws.Rows.Item(8).AutoFit()

# SmartScore text cells (kept as text, matching the original text-formatted values)
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "0.575"
$ws.Cells.Item(8,7).ClearFormats()
$ws.Cells.Item(8,10).NumberFormat = "@"
$ws.Cells.Item(8,10).Value = "0.510"
$ws.Cells.Item(8,10).ClearFormats()
$ws.Cells.Item(8,13).NumberFormat = "@"
$ws.Cells.Item(8,13).Value = "0.509"
$ws.Cells.Item(8,13).ClearFormats()
$ws.Cells.Item(8,16).NumberFormat = "@"
$ws.Cells.Item(8,16).Value = "0.650"
$ws.Cells.Item(8,16).ClearFormats()
$ws.Cells.Item(8,19).NumberFormat = "@"
$ws.Cells.Item(8,19).Value = "0.587"
$ws.Cells.Item(8,19).ClearFormats()
$ws.Cells.Item(8,22).NumberFormat = "@"
$ws.Cells.Item(8,22).Value = "0.552"
$ws.Cells.Item(8,22).ClearFormats()
$ws.Cells.Item(8,25).NumberFormat = "@"
$ws.Cells.Item(8,25).Value = "0.664"
$ws.Cells.Item(8,25).ClearFormats()
$ws.Cells.Item(8,28).NumberFormat = "@"
$ws.Cells.Item(8,28).Value = "0.589"
$ws.Cells.Item(8,28).ClearFormats()
$ws.Cells.Item(8,31).NumberFormat = "@"
$ws.Cells.Item(8,31).Value = "0.576"
$ws.Cells.Item(8,31).ClearFormats()

# --- Step 2: Convert row 7 SmartScore cells from text to real numeric values ---
$ws.Cells.Item(7,7).Value = 0.575
$ws.Cells.Item(7,10).Value = 0.51
$ws.Cells.Item(7,13).Value = 0.509
$ws.Cells.Item(7,16).Value = 0.65
$ws.Cells.Item(7,19).Value = 0.587
$ws.Cells.Item(7,22).Value = 0.552
$ws.Cells.Item(7,25).Value = 0.664
$ws.Cells.Item(7,28).Value = 0.589
$ws.Cells.Item(7,31).Value = 0.576

Write-Output "done"
